$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update price and quantity values for row 2
$ws.Range("C2").Value = 25000
$ws.Range("D2").Value = 25

# Move the active selection to F6
$ws.Range("F6").Select()
